# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2024-09-04 (serial 45539) to 2024-09-05 (serial 45540).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45539) {
        $cell.Value2 = 45540
    }
}
